$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells stay as text (matches inlineStr in source)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.759.20"
$ws.Range("E2").Value = "  +4.86%  "
$ws.Range("D3").Value = "2.275.67"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "231.26"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "63.80"
$ws.Range("E7").Value = "  +5.96%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.429"
$ws.Range("E9").Value = "  +6.52%  "
$ws.Range("D10").Value = "0.105"
$ws.Range("E10").Value = "  +16.52%  "
$ws.Range("D11").Value = "57.32"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "25.84"
$ws.Range("E12").Value = "  +14.74%  "
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "2.610.24"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "5.88"
$ws.Range("E16").Value = "  +4.48%  "
$ws.Range("D17").Value = "0.818"
$ws.Range("E17").Value = "  +2.27%  "
$ws.Range("D18").Value = "2.287.13"
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("D19").Value = "43.643.12"
$ws.Range("E19").Value = "  +4.71%  "
$ws.Range("D20").Value = "0.0000102"
$ws.Range("E20").Value = "  +11.73%  "
$ws.Range("D21").Value = "73.35"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "249.02"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("E25").Value = "  +5.30%  "
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("D27").Value = "9.83"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "171.85"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").Value = "21.13"
$ws.Range("E29").Value = "  +6.17%  "
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").Value = "2.78"
$ws.Range("E32").Value = "  +8.42%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "0.0689"
$ws.Range("E34").Value = "  +5.14%  "
$ws.Range("D35").Value = "5.08"
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("D36").Value = "4.70"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("E37").Value = "  +6.05%  "
$ws.Range("E38").Value = "  +2.93%  "
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "8.37"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("D43").Value = "10.50"
$ws.Range("E43").Value = "  +19.54%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0961"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "17.14"
$ws.Range("E45").Value = "  +3.19%  "
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "97.32"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "4.39"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Value = "1.476.47"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "2.32"
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").Value = "0.000204"
$ws.Range("E51").Value = "  -15.41%  "
